$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores numeric-looking strings (e.g. "70.118.94",
# "0.998") as literal TEXT in the source workbook. Assigning such a string
# via .Value would normally let Excel auto-convert it to a number/date, so
# we first mark each touched D-column cell as Text ("@") before writing its
# new value, keeping the cell type in sync with the original inlineStr data.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "70.118.94"
$ws.Range("E2").Value = "  +3.00%  "
$ws.Range("D3").Value = "2.581.79"
$ws.Range("E3").Value = "  +2.54%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("D5").Value = "603.38"
$ws.Range("E5").Value = "  +2.28%  "
$ws.Range("D6").Value = "178.78"
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").Value = "0.521"
$ws.Range("E8").Value = "  +0.92%  "
$ws.Range("D9").Value = "2.581.63"
$ws.Range("E9").Value = "  +2.56%  "
$ws.Range("D10").Value = "0.161"
$ws.Range("E10").Value = "  +13.80%  "
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").Value = "0.347"
$ws.Range("E12").Value = "  +1.63%  "
$ws.Range("D13").Value = "5.05"
$ws.Range("E13").Value = "  +1.61%  "
$ws.Range("D14").Value = "0.0000184"
$ws.Range("E14").Value = "  +6.00%  "
$ws.Range("D15").Value = "26.45"
$ws.Range("E15").Value = "  +2.37%  "
$ws.Range("D16").Value = "2.982.53"
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("D17").Value = "69.894.30"
$ws.Range("E17").Value = "  +2.95%  "
$ws.Range("D18").Value = "2.541.06"
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("D19").Value = "7.76"
$ws.Range("E19").Value = "  +2.27%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "11.22"
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "365.87"
$ws.Range("E21").Value = "  +3.51%  "
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").Value = "71.05"
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("D25").Value = "4.33"
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("B26").Value = "SuiNetwork"
$ws.Range("C26").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D26").Value = "1.75"
$ws.Range("E26").Value = "  -2.09%  "
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").Value = "9.26"
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("D28").Value = "0.988"
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("D29").Value = "0.0₃0942"
$ws.Range("E29").Value = "  +1.91%  "
$ws.Range("D30").Value = "523.16"
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("D31").Value = "7.85"
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("D32").Value = "1.29"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").Value = "1.81"
$ws.Range("E33").Value = "  +1.70%  "
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "0.121"
$ws.Range("E35").Value = "  -1.87%  "
$ws.Range("D36").Value = "161.60"
$ws.Range("E36").Value = "  -1.19%  "
$ws.Range("D37").Value = "19.01"
$ws.Range("E37").Value = "  +3.04%  "
$ws.Range("D38").Value = "18.93"
$ws.Range("E38").Value = "  +1.39%  "
$ws.Range("D39").Value = "1.36"
$ws.Range("E39").Value = "  +0.74%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "1.78"
$ws.Range("E40").Value = "  +1.43%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "4.98"
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "2.52"
$ws.Range("E43").Value = "  +1.21%  "
$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").Value = "0.328"
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("D45").Value = "39.08"
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("D46").Value = "153.08"
$ws.Range("E46").Value = "  +4.12%  "
$ws.Range("D47").Value = "3.66"
$ws.Range("E47").Value = "  +2.76%  "
$ws.Range("D48").Value = "0.529"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("D49").Value = "0.0₆0263"
$ws.Range("E49").Value = "  +0.72%  "
$ws.Range("D50").Value = "1.63"
$ws.Range("E50").Value = "  +1.65%  "
$ws.Range("D51").Value = "0.0745"
$ws.Range("E51").Value = "  +0.15%  "
